$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.825
$ws.Range("D2").Value = 0.263
$ws.Range("E2").Value = 0.62

$ws.Range("B3").Value = 4.692
$ws.Range("D3").Value = 1.498
$ws.Range("E3").Value = 0.252

$ws.Range("B4").Value = 4524.858
$ws.Range("D4").Value = 1444.132
$ws.Range("E4").Value = [double]"3.003321566616797E-11"

$ws.Range("B5").Value = 44046.367
$ws.Range("D5").Value = 14057.622
$ws.Range("E5").Value = [double]"1.09680721783061E-15"

$ws.Range("B6").Value = 5105.875
$ws.Range("D6").Value = 1629.566
$ws.Range("E6").Value = [double]"1.748885486397493E-11"

$ws.Range("B7").Value = 0.01237164500000661
$ws.Range("D7").Value = 0.003948473566969992
$ws.Range("E7").Value = 0.951

$ws.Range("B8").Value = 75.857
$ws.Range("D8").Value = 24.21
$ws.Range("E8").Value = 0.000824185498638506

$ws.Range("B9").Value = 0.01237164500000884
$ws.Range("D9").Value = 0.003948473566970702
$ws.Range("E9").Value = 0.951

$ws.Range("B10").Value = 1.13
$ws.Range("D10").Value = 0.361
$ws.Range("E10").Value = 0.563

$ws.Range("B11").Value = 1.661
$ws.Range("D11").Value = 0.53
$ws.Range("E11").Value = 0.485

$ws.Range("B12").Value = 5.378
$ws.Range("D12").Value = 1.716
$ws.Range("E12").Value = 0.223

$ws.Range("B13").Value = 4.692
$ws.Range("D13").Value = 1.498
$ws.Range("E13").Value = 0.252

$ws.Range("B14").Value = 7.265
$ws.Range("D14").Value = 2.319
$ws.Range("E14").Value = 0.162

$ws.Range("B15").Value = 0.714
$ws.Range("D15").Value = 0.228
$ws.Range("E15").Value = 0.644

$ws.Range("B16").Value = 4.481
$ws.Range("D16").Value = 1.43
$ws.Range("E16").Value = 0.262

$ws.Range("B17").Value = 5.59
$ws.Range("D17").Value = 1.784
$ws.Range("E17").Value = 0.214

$ws.Range("B18").Value = 0.825
$ws.Range("D18").Value = 0.263
$ws.Range("E18").Value = 0.62

$ws.Range("B19").Value = 4.69
$ws.Range("D19").Value = 1.497
$ws.Range("E19").Value = 0.252

$ws.Range("B20").Value = 4526.381
$ws.Range("D20").Value = 1444.617
$ws.Range("E20").Value = [double]"2.998804759619007E-11"

$ws.Range("B21").Value = 54.08
$ws.Range("D21").Value = 17.26
$ws.Range("E21").Value = 0.002467793261966397

$ws.Range("B22").Value = 6.529
$ws.Range("D22").Value = 2.084
$ws.Range("E22").Value = 0.183

$ws.Range("B23").Value = 5025.071
$ws.Range("D23").Value = 1603.777
$ws.Range("E23").Value = [double]"1.878365315325018E-11"

$ws.Range("B24").Value = 5025.071
$ws.Range("D24").Value = 1603.777
$ws.Range("E24").Value = [double]"1.878365315325036E-11"

$ws.Range("B25").Value = 0.435
$ws.Range("D25").Value = 0.139
$ws.Range("E25").Value = 0.718

$ws.Range("B26").Value = 54.08
$ws.Range("D26").Value = 17.26
$ws.Range("E26").Value = 0.002467793261966397

$ws.Range("B27").Value = 0.11
$ws.Range("D27").Value = 0.03499613678177024
$ws.Range("E27").Value = 0.856

$ws.Range("B28").Value = 4526.381
$ws.Range("D28").Value = 1444.617
$ws.Range("E28").Value = [double]"2.998804759619067E-11"

$ws.Range("B29").Value = 2.057
$ws.Range("D29").Value = 0.656
$ws.Range("E29").Value = 0.439

$ws.Range("B30").Value = 28.199
$ws.Range("C30").Value = 9
